$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions (K1, L1) ---
# shared-string insertion order matters for index parity with target: "test" (16) then "comments" (17)
$ws.Range("L1").Value = "test"
$ws.Range("K1").Value = "comments"

# --- Row 4 value updates ---
$ws.Range("E4").Value = 35
$ws.Range("F4").Value = 0.1969
$ws.Range("G4").Value = 0.991
$ws.Range("H4").Value = 0.99
$ws.Range("I4").Value = 0.995
$ws.Range("J4").Value = 0.935

# Remove the old K4 "RISCRIVI" comment cell entirely
$ws.Range("K4").ClearContents()

# --- New comment values in column L for a few existing rows ---
$ws.Range("L5").Value = 0.9379
$ws.Range("L8").Value = 0.9445

# --- New rows 12 and 13 ---
$ws.Range("A12").Value = "11s"
$ws.Range("B12").Value = 512
$ws.Range("C12").Value = "no"
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 0.39
$ws.Range("G12").Value = 0.998
$ws.Range("H12").Value = 0.997
$ws.Range("I12").Value = 0.995
$ws.Range("J12").Value = 0.965
$ws.Range("L12").Value = 0.9434

$ws.Range("A13").Value = "11s"
$ws.Range("B13").Value = 640
$ws.Range("C13").Value = "light"
$ws.Range("E13").Value = 249
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 0.997
$ws.Range("H13").Value = 0.998
$ws.Range("I13").Value = 0.995
$ws.Range("J13").Value = 0.974
$ws.Range("L13").Value = 0.9369

# --- Fill D3:D13 with the (F*3600)/E formula (drag-fill style) ---
$ws.Range("D3:D13").FormulaR1C1 = "=(RC[2]*3600)/RC[1]"

# --- Apply a fill-related style touch to columns A:B for all data rows (2-13) ---
# (mirrors the workbook author selecting A2:B13 and toggling a fill swatch)
$fmtRange = $ws.Range("A2:B13")
$fmtRange.Interior.Pattern = 1
$fmtRange.Interior.ColorIndex = -4142

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 145
$ws.Range("G17").Select()
